$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (8-15): from_bus (C), to_bus (D), in_service (E) ---

# Row 8 - extr1
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

# Row 9 - extr2
$ws.Range("C9").Value = 16
$ws.Range("E9").Value = $true

# Row 10 - extr3
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

# Row 11 - extr4
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9

# Row 12 - extr5
$ws.Range("C12").Value = 10

# Row 13 - extr6
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $true

# Row 14 - extr7
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11

# Row 15 - extr8
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11

# --- Add two new rows (16-17) for line7 and line8 ---

# Column A on the new rows uses the same style as the rest of column A
# (bold, thin-bordered, centered cell format, style index 1). Copy that
# formatting from the row above instead of re-deriving it property by
# property, so the workbook ends up re-using the existing cell style
# (xf index) instead of synthesizing new, unused ones.
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A17").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Row 16 - line7
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "line7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $false

# Row 17 - line8
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "line8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $false
